# Updates odds values in Jogos_da_Semana_FlashScore_2025-02-06 sheet
# to match the refreshed FlashScore odds data (per commit "Atualizando o arquivo XLSX").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("G3").Value = 2.5
$ws.Range("I3").Value = 3.3
$ws.Range("Y3").Value = 1.63
$ws.Range("AC3").Value = 5.5
$ws.Range("AF3").Value = 23
$ws.Range("AJ3").Value = 6
$ws.Range("AO3").Value = 15
# Row 4
$ws.Range("Y4").Value = 1.69
# Row 5
$ws.Range("S5").Value = 3.4
$ws.Range("T5").Value = 1.33
$ws.Range("Y5").Value = 1.73
$ws.Range("Z5").Value = 2.08
# Row 6
$ws.Range("G6").Value = 2.4
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2.88
$ws.Range("J6").Value = 3
$ws.Range("L6").Value = 3.4
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 12
$ws.Range("S6").Value = 1.83
$ws.Range("T6").Value = 2.03
$ws.Range("Y6").Value = 1.33
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 9.5
$ws.Range("AF6").Value = 23
$ws.Range("AG6").Value = 19
$ws.Range("AH6").Value = 26
$ws.Range("AJ6").Value = 6.5
$ws.Range("AN6").Value = 10
$ws.Range("AR6").Value = 21
# Row 7
$ws.Range("G7").Value = 1.72
$ws.Range("H7").Value = 3.25
$ws.Range("I7").Value = 5.3
$ws.Range("J7").Value = 2.25
$ws.Range("K7").Value = 2.07
$ws.Range("L7").Value = 5.5
$ws.Range("M7").Value = 1.09
$ws.Range("N7").Value = 6.2
$ws.Range("P7").Value = 2.77
$ws.Range("S7").Value = 2.15
$ws.Range("T7").Value = 1.62
$ws.Range("W7").Value = 3.7
$ws.Range("AA7").Value = 2
$ws.Range("AB7").Value = 1.72
$ws.Range("AD7").Value = 7.3
$ws.Range("AF7").Value = 13.5
$ws.Range("AI7").Value = 6.2
$ws.Range("AJ7").Value = 6.4
$ws.Range("AK7").Value = 17.5
$ws.Range("AL7").Value = 100
$ws.Range("AO7").Value = 30
$ws.Range("AQ7").Value = 110
$ws.Range("AR7").Value = 65
# Row 9
$ws.Range("I9").Value = 9
$ws.Range("J9").Value = 1.88
$ws.Range("K9").Value = 2.2
$ws.Range("L9").Value = 7.8
$ws.Range("N9").Value = 6.7
$ws.Range("S9").Value = 2.02
$ws.Range("X9").Value = 1.27
$ws.Range("Z9").Value = 2.7
$ws.Range("AA9").Value = 2.3
$ws.Range("AB9").Value = 1.55
$ws.Range("AC9").Value = 5.2
$ws.Range("AD9").Value = 5.5
$ws.Range("AF9").Value = 8.5
$ws.Range("AI9").Value = 6.7
$ws.Range("AN9").Value = 19
$ws.Range("AO9").Value = 65
$ws.Range("AP9").Value = 28
$ws.Range("AQ9").Value = 300
# Row 10
$ws.Range("Y10").Value = 1.25
$ws.Range("AA10").Value = 1.54
# Row 11
$ws.Range("G11").Value = 1.95
$ws.Range("H11").Value = 3.5
$ws.Range("I11").Value = 3.75
$ws.Range("J11").Value = 2.6
$ws.Range("K11").Value = 2.2
$ws.Range("L11").Value = 4
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 13
$ws.Range("S11").Value = 1.83
$ws.Range("T11").Value = 2.03
$ws.Range("AD11").Value = 10
$ws.Range("AM11").Value = 201
$ws.Range("AN11").Value = 12
$ws.Range("AO11").Value = 19
$ws.Range("AR11").Value = 29
# Row 12
$ws.Range("H12").Value = 6
$ws.Range("I12").Value = 10
$ws.Range("J12").Value = 1.62
$ws.Range("K12").Value = 2.63
$ws.Range("M12").Value = 1.01
$ws.Range("N12").Value = 17
$ws.Range("O12").Value = 1.14
$ws.Range("P12").Value = 5
$ws.Range("S12").Value = 1.5
$ws.Range("T12").Value = 2.5
$ws.Range("U12").Value = 1.78
$ws.Range("V12").Value = 2.03
$ws.Range("W12").Value = 2.2
$ws.Range("X12").Value = 1.62
$ws.Range("Y12").Value = 1.25
$ws.Range("Z12").Value = 3.75
$ws.Range("AC12").Value = 8.5
$ws.Range("AE12").Value = 9.5
$ws.Range("AG12").Value = 11
$ws.Range("AI12").Value = 17
$ws.Range("AK12").Value = 23
$ws.Range("AN12").Value = 26
$ws.Range("AP12").Value = 29
# Row 13
$ws.Range("G13").Value = 1.65
$ws.Range("H13").Value = 4.2
$ws.Range("K13").Value = 2.4
$ws.Range("L13").Value = 4.33
$ws.Range("N13").Value = 10.5
$ws.Range("O13").Value = 1.17
$ws.Range("P13").Value = 4.5
$ws.Range("S13").Value = 1.57
$ws.Range("T13").Value = 2.35
$ws.Range("U13").Value = 1.9
$ws.Range("V13").Value = 1.9
$ws.Range("W13").Value = 2.38
$ws.Range("X13").Value = 1.53
$ws.Range("Y13").Value = 1.29
$ws.Range("Z13").Value = 3.5
$ws.Range("AA13").Value = 1.62
$ws.Range("AB13").Value = 2.2
$ws.Range("AC13").Value = 9.5
$ws.Range("AD13").Value = 9.5
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 17
$ws.Range("AJ13").Value = 8.5
$ws.Range("AN13").Value = 17
$ws.Range("AO13").Value = 26
$ws.Range("AR13").Value = 29
# Row 14
$ws.Range("M14").Value = 1.04
$ws.Range("N14").Value = 9
$ws.Range("O14").Value = 1.25
$ws.Range("P14").Value = 3.75
$ws.Range("S14").Value = 1.83
$ws.Range("T14").Value = 1.98
$ws.Range("W14").Value = 3
$ws.Range("X14").Value = 1.36
# Row 15
$ws.Range("G15").Value = 1.65
$ws.Range("H15").Value = 4.33
$ws.Range("AA15").Value = 1.5
$ws.Range("AB15").Value = 2.5
$ws.Range("AC15").Value = 11
$ws.Range("AF15").Value = 13
$ws.Range("AJ15").Value = 9
$ws.Range("AM15").Value = 101
$ws.Range("AN15").Value = 19
$ws.Range("AO15").Value = 29
$ws.Range("AR15").Value = 34
# Row 16
$ws.Range("N16").Value = 13
# Row 18
$ws.Range("G18").Value = 1.98
$ws.Range("H18").Value = 3.35
$ws.Range("I18").Value = 3.45
$ws.Range("J18").Value = 2.57
$ws.Range("L18").Value = 4
$ws.Range("O18").Value = 1.35
$ws.Range("P18").Value = 2.7
$ws.Range("S18").Value = 2.02
$ws.Range("T18").Value = 1.62
$ws.Range("W18").Value = 3.35
$ws.Range("X18").Value = 1.23
$ws.Range("Y18").Value = 1.42
$ws.Range("Z18").Value = 2.47
$ws.Range("AC18").Value = 6.5
$ws.Range("AD18").Value = 8.75
$ws.Range("AE18").Value = 8.75
$ws.Range("AF18").Value = 17
$ws.Range("AG18").Value = 17.5
$ws.Range("AH18").Value = 32
$ws.Range("AI18").Value = 8.5
$ws.Range("AJ18").Value = 6.5
$ws.Range("AK18").Value = 17
$ws.Range("AN18").Value = 9
$ws.Range("AO18").Value = 17.5
$ws.Range("AP18").Value = 12.5
$ws.Range("AQ18").Value = 50
$ws.Range("AR18").Value = 35
